# Apply updated voltage-magnitude (vm_pu) results for the 380 kV case.
# Rows 2-25 (bus indices 0-23) get new values in columns B,C,D,E,F,I,J,K,L,M,N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (in the order the data values appear in each row below)
$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

$data = @(
    @(2, 1.02, 1.02447107144795, 1.03142051084628, 1.024933484179309, 1.022953973049677, 1.033668293642911, 1.029645501202792, 1.034228679803083, 1.027760533228763, 1.025786835219658, 1.01385080037022),
    @(3, 1.02, 1.025701623551068, 1.031995083186469, 1.025985973084647, 1.024825367792706, 1.033947218055538, 1.03051352114588, 1.03461243975731, 1.028619551724816, 1.027462103305554, 1.01414439183759),
    @(4, 1.02, 1.026496591638564, 1.032366315163305, 1.026666233498146, 1.026034688194842, 1.034125765638977, 1.031073450099932, 1.034859503943135, 1.02917403378132, 1.028544118149551, 1.014333618873426),
    @(5, 1.02, 1.026830494331151, 1.0325222482718, 1.026952033929406, 1.026542716971301, 1.034200364694013, 1.031308431735906, 1.034963069918667, 1.029406815572278, 1.028998530876138, 1.014412992359335),
    @(6, 1.02, 1.026886540500582, 1.032548422337958, 1.027000010559862, 1.026627995940887, 1.034212863118062, 1.031347862101093, 1.034980441533296, 1.029445881842342, 1.029074801758704, 1.014426309136439),
    @(7, 1.02, 1.026501054442895, 1.032368399271064, 1.026670053084756, 1.026041477936364, 1.03412676425079, 1.031076591552486, 1.034860888974301, 1.02917714548614, 1.02855019185285, 1.014334680162171),
    @(8, 1.02, 1.024887210564433, 1.031614805125333, 1.025289339257719, 1.023586758111499, 1.033762958708242, 1.029939214634781, 1.034358633230645, 1.028051126158818, 1.026353420030517, 1.013950176167713),
    @(9, 1.02, 1.02203337743018, 1.03028262605436, 1.022850293673801, 1.019248418977741, 1.033107023180064, 1.027921538457864, 1.033463968140149, 1.026056372749567, 1.022466628706994, 1.013266865649747),
    @(10, 1.02, 1.020123765121997, 1.029391643866302, 1.021219979060958, 1.016346771854904, 1.032659677651599, 1.026567133676163, 1.032861020673652, 1.02471922960314, 1.019864090712372, 1.012807378268729),
    @(11, 1.02, 1.019295138101574, 1.029005157377736, 1.020512969751384, 1.015087911529232, 1.032463572783711, 1.025978406808027, 1.032598387161516, 1.02413845210242, 1.018734311259047, 1.012607462886422),
    @(12, 1.02, 1.018987079733857, 1.028861495810833, 1.020250189593162, 1.014619934364593, 1.032390368616288, 1.025759383232352, 1.032500599244317, 1.02392245309654, 1.018314216728585, 1.012533060609215),
    @(13, 1.02, 1.019053171589395, 1.028892316383419, 1.020306564384579, 1.014720334520566, 1.032406087565333, 1.025806380182569, 1.032521585692353, 1.023968798008758, 1.018404348754992, 1.012549026725979),
    @(14, 1.02, 1.019269679432127, 1.02899328437653, 1.020491251661206, 1.015049236241461, 1.032457529098583, 1.025960309299167, 1.032590308762655, 1.02412060313763, 1.018699595255574, 1.012601315733234),
    @(15, 1.02, 1.019403041177354, 1.029055480394417, 1.020605021537963, 1.015251832572096, 1.032489175905027, 1.026055104360484, 1.032632620265459, 1.024214099070867, 1.018881447096683, 1.01263351348556),
    @(16, 1.02, 1.020178720927146, 1.029417279236453, 1.021266877980345, 1.016430265784678, 1.032672641780984, 1.02660615752138, 1.03287841800346, 1.024757735921394, 1.019939008938529, 1.012820625766417),
    @(17, 1.02, 1.020664810888749, 1.029644042396791, 1.021681752915892, 1.017168804788954, 1.032787081255694, 1.026951210159254, 1.033032184033222, 1.02509826399079, 1.020601613134481, 1.012937739856482),
    @(18, 1.02, 1.02094817028703, 1.029776243411735, 1.021923639277806, 1.017599348892629, 1.032853600264978, 1.027152255715688, 1.033121723341958, 1.025296716238824, 1.020987823631193, 1.013005958543906),
    @(19, 1.02, 1.021044760097263, 1.02982130933852, 1.02200609890663, 1.017746114220706, 1.032876242284154, 1.02722077024781, 1.033152228543966, 1.025364354273629, 1.021119465151074, 1.013029203794849),
    @(20, 1.02, 1.020612675526582, 1.029619719695577, 1.021637251465309, 1.017089590780157, 1.032774826948616, 1.026914211839992, 1.033015701895577, 1.025061746374606, 1.02053055055622, 1.012925184147613),
    @(21, 1.02, 1.019205930755014, 1.028963554672719, 1.020436870462776, 1.014952393509004, 1.032442390855281, 1.025914990549517, 1.03257007800098, 1.024075907863087, 1.018612664893065, 1.012585921932686),
    @(22, 1.02, 1.018319892751945, 1.028550399977765, 1.019681184037032, 1.013606444111978, 1.032231279492179, 1.025284746713374, 1.032288541528139, 1.02345449407613, 1.017404237407393, 1.012371775808616),
    @(23, 1.02, 1.018789748358567, 1.028769477865399, 1.020081880087541, 1.014320171956746, 1.032343392703938, 1.025619041416105, 1.032437918045458, 1.023784068237545, 1.018045096609844, 1.01248537868207),
    @(24, 1.02, 1.020636233773258, 1.029630710281632, 1.02165736007378, 1.017125384900007, 1.032780364858041, 1.026930930460541, 1.033023149927026, 1.025078247644677, 1.020562661525793, 1.012930857816101),
    @(25, 1.02, 1.022772382587767, 1.030627530378984, 1.023481584596491, 1.020371584454204, 1.033278365943057, 1.028444776133628, 1.033696404487448, 1.026573337144058, 1.023473403583669, 1.013444208253709)
)

foreach ($row in $data) {
    $r = $row[0]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $c = $cols[$i]
        $val = $row[$i + 1]
        $ws.Cells.Item($r, $c).Value = $val
    }
}
